$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# The "TabStats_V02" timing table (A11:G18) had its Timing_Exec column (F)
# values updated for the first four data rows (F2 threads no longer use the
# combin list - that logic moved to a different file, hence new timings).
# These are stored as text in the workbook. "0.724 sec" already reads as
# text (it isn't a pure number), so a plain Value assignment keeps it text.
$ws.Range("F12").Value = "0.724 sec"

# The remaining new values ("2.473", "12.469", "36.197") look like plain
# numbers, so a plain Value assignment would store them as numeric cells.
# Force text format while writing them, then clear the (now unneeded)
# number format again so the cell formatting is left as it was - only the
# text content changes, same as the original file.
function Set-TimingText($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TimingText "F13" "2.473"
Set-TimingText "F15" "36.197"
Set-TimingText "F14" "12.469"
